$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 510, shifting the existing data (rows 510:526) down to 512:528.
$ws.Rows("510:511").Insert()

# New row 510 — Mercado Mayorista Lo Valledor de Santiago, Poroto verde, Magnum/Primera, Peru.
$ws.Range("A510").Value = 6
$ws.Range("B510").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C510").Value = "Metropolitana"
$ws.Range("D510").Value = 44509
$ws.Range("E510").Value = 13
$ws.Range("F510").Value = 100112031
$ws.Range("G510").Value = "Poroto verde"
$ws.Range("H510").Value = "Magnum"
$ws.Range("I510").Value = "Primera"
$ws.Range("J510").Value = 350
$ws.Range("K510").Value = 30000
$ws.Range("L510").Value = 35000
$ws.Range("M510").Value = 31714
$ws.Range("N510").Value = "$/malla 25 kilos"
$ws.Range("O510").Value = "Perú"
$ws.Range("P510").Value = 1269
$ws.Range("Q510").Value = 25
$ws.Range("R510").Value = "Hortaliza"

# New row 511 — same market/date, Magnum/Segunda, Peru.
$ws.Range("A511").Value = 6
$ws.Range("B511").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C511").Value = "Metropolitana"
$ws.Range("D511").Value = 44509
$ws.Range("E511").Value = 13
$ws.Range("F511").Value = 100112031
$ws.Range("G511").Value = "Poroto verde"
$ws.Range("H511").Value = "Magnum"
$ws.Range("I511").Value = "Segunda"
$ws.Range("J511").Value = 100
$ws.Range("K511").Value = 25000
$ws.Range("L511").Value = 25000
$ws.Range("M511").Value = 25000
$ws.Range("N511").Value = "$/malla 25 kilos"
$ws.Range("O511").Value = "Perú"
$ws.Range("P511").Value = 1000
$ws.Range("Q511").Value = 25
$ws.Range("R511").Value = "Hortaliza"
